# Add two new rows to the "tableau de bord" features table, matching the
# formatting Word itself clones from the table's last row (cell widths,
# vAlign/jc, the 22-half-point run size on the feature-name cell), then
# mark the "fichier CSV" row's checkmark cell with the _GoBack bookmark
# (moving it off the trailing paragraph, which becomes empty).

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# --- Row: "Gestion de l'affichage des messages d'alerte" -> X in column 4 ---
$row1 = $table.Rows.Add()
$row1.Cells.Item(1).Range.Text = "Gestion de l’affichage des messages d’alerte"
$row1.Cells.Item(4).Range.Text = "X"

# --- Row: "Fonctionnalité de téléchargement ..." -> X in column 3 ---
$row2 = $table.Rows.Add()
$row2.Cells.Item(1).Range.Text = "Fonctionnalité de téléchargement du fichier CSV de la liste des inscrits à une épreuve"

# Write a temporary trailing character after the "X" so the insertion
# point sits right after a real character (not on the cell's end-of-cell
# mark), add the _GoBack bookmark there (which, being a document-unique
# name, relocates it off the trailing paragraph), then delete the
# temporary character again. The bookmark stays anchored right after "X".
$cell3 = $row2.Cells.Item(3)
$cell3.Range.Text = "XZ"
$afterX = $cell3.Range.Start + 1
$bmRange = $d.Range($afterX, $afterX)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($afterX, $afterX + 1).Delete()
